$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, matching style/formatting of existing header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill in Save column values for existing data rows
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
